$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 16
$ws.Range("I2").Value = 0.01711229946524064
$ws.Range("J2").Value = 0.8
$ws.Range("K2").Value = 0.03350785340314136
